$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column numbers in the sheet BEFORE the edit:
#   DR = 122  (last existing price-snapshot column)
#   DS = 123  ("nom" column)
#   DT = 124  ("url_produit" column)
$colLastPrice = 122
$colNewSnapshot = 123

# Insert a brand-new column right before the existing "nom" column (DS).
# This shifts "nom" from DS -> DT and "url_produit" from DT -> DU, one
# column to the right, exactly like Excel's own Insert Column command.
$ws.Columns("DS:DS").Insert()

# Header of the freshly inserted column: timestamp of this price check.
$ws.Range("DS1").Value = "2026-02-02 09:35:37"

# For every product row, if the previous run recorded a price in DR,
# the price is unchanged since that check, so carry it over into the
# new DS column. Rows with no price in DR (out of stock) stay blank.
$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $price = $ws.Cells.Item($r, $colLastPrice).Value()
    if ($price -ne $null -and $price -ne "") {
        $ws.Cells.Item($r, $colNewSnapshot).Value = $price
    }
}
